$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (issue number and week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# --- Crime statistics table updates ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -100
$ws.Range("L15").Value = 175
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 43
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = -30.645161290322
$ws.Range("L16").Value = -28.333333333333
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 31.25
$ws.Range("I17").Value = 183
$ws.Range("J17").Value = 157
$ws.Range("K17").Value = 16.56050955414
$ws.Range("L17").Value = 7.647058823529
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 37
$ws.Range("K18").Value = -28.846153846153
$ws.Range("L18").Value = -31.481481481481
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 151
$ws.Range("J19").Value = 187
$ws.Range("K19").Value = -19.251336898395
$ws.Range("L19").Value = -30.092592592592
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 6
$ws.Range("H20").Value = -62.5
$ws.Range("I20").Value = 116
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = -3.333333333333
$ws.Range("L20").Value = -17.142857142857
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -22.222222222222
$ws.Range("F21").Value = 51
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = -21.538461538461
$ws.Range("I21").Value = 552
$ws.Range("J21").Value = 588
$ws.Range("K21").Value = -6.122448979591
$ws.Range("L21").Value = -15.076923076923
$ws.Range("C24").Value = 6
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 32
$ws.Range("H24").Value = 159.375
$ws.Range("I24").Value = 407
$ws.Range("J24").Value = 387
$ws.Range("K24").Value = 5.16795865633
$ws.Range("L24").Value = -3.095238095238
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = -10
$ws.Range("I25").Value = 72
$ws.Range("J25").Value = 110
$ws.Range("K25").Value = -34.545454545454
$ws.Range("L25").Value = -36.283185840708
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 41
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = 13.888888888888
$ws.Range("I26").Value = 329
$ws.Range("J26").Value = 312
$ws.Range("K26").Value = 5.448717948717
$ws.Range("L26").Value = 31.075697211155
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -100
$ws.Range("L27").Value = 80
$ws.Range("G28").Value = 5
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = -75
$ws.Range("L29").Value = -50
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = -66.666666666666
$ws.Range("L30").Value = -33.333333333333

# --- Cells switching from numeric to text placeholder "0" (style 13 / shared string "0") ---
# Source of target format: D15 (style 13, text "0")
$ws.Range("F15").Value = "'0"
$ws.Range("C16").Value = "'0"
$ws.Range("D27").Value = "'0"
$ws.Range("F27").Value = "'0"
$ws.Range("D28").Value = "'0"
$ws.Range("D15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("D28").PasteSpecial(-4122)

# --- Cells switching from numeric to text placeholder "***.*" (style 13 / shared string "***.*") ---
# Source of target format: E15 (style 13, text "***.*")
$ws.Range("E27").Value = "'***.*"
$ws.Range("E28").Value = "'***.*"
$ws.Range("E15").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E28").PasteSpecial(-4122)

# --- Cells switching from text placeholder to numeric value (style 14 / #,##0 number format) ---
# Source of target format: G15 (style 14, numeric)
$ws.Range("C18").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("G15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("F30").PasteSpecial(-4122)

$excel.CutCopyMode = $false
